$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(39400,0.2592894031801904),
    @(39583,1.528044834838681),
    @(39765,-0.4231639722327287),
    @(39948,-3.800765243909694),
    @(40130,0.4426058137615883),
    @(40310,0.1806671024190649),
    @(40494,2.3),
    @(40676,0.3647642815062397),
    @(40862,0.5),
    @(41044,0.5),
    @(41228,0.2),
    @(41409,0.06328053690715763),
    @(41592,0.7),
    @(41774,0.8179240379891439),
    @(41957,-0.08509419607057112),
    @(42137,0.699690853462215),
    @(42321,0.4386429368392442),
    @(42503,0.268690731786478),
    @(42689,0.1911347634769811),
    @(42867,0.4177287092911968),
    @(43053,0.633202420130246),
    @(43145,0.735487593389081),
    @(43235,0.2972792544561997),
    @(43326,0.366083231974244),
    @(43418,0.5),
    @(43510,-0.1988240152552123),
    @(43600,0.4),
    @(43691,0.3842522401052975),
    @(43783,-0.2425341822259668),
    @(43875,0.1963547214277668),
    @(43966,-2.2),
    @(44068,-9.690955511628317),
    @(44159,8.480305474188015),
    @(44251,0.34062844669036),
    @(44341,-1.80694255667963),
    @(44432,1.633091849519602),
    @(44525,1.694275412034003),
    @(44617,-0.3471888372093019),
    @(44706,0.2354020960014367),
    @(44798,0.1397248471449046),
    @(44890,0.4001468552058896),
    @(44981,-0.4353111306276531),
    @(45071,-0.3352883615304449),
    @(45163,0.01856969850348378),
    @(45254,-0.1296233709318813),
    @(45345,-0.2870570324401456),
    @(45436,0.2135990538361767),
    @(45534,-0.06675723343442996),
    @(45618,0.1052088781249552),
    @(45713,-0.2005382402049349),
    @(45800,0.4116519254275204),
    @(45891,-0.2766843069119602)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -gt 53) {
    $ws.Range("A54:B" + $lastRow).EntireRow.Delete()
}
